$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C - copy formatting (bold/border/centered) from B1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "bry"

# Updated B column values + new C column values
$data = @(
    @(2,  6.000302,           5.165224),
    @(3,  4.23926,            3.404182),
    @(4,  4.424898999999999,  3.589820999999999),
    @(5,  4.814907000000001,  3.979828999999999),
    @(6,  9.395856,           8.560777999999997),
    @(7,  9.377315000000001,  8.542236999999997),
    @(8,  8.202413999999997,  7.367335999999997),
    @(9,  6.632353999999999,  5.797275999999999),
    @(10, 4.819519000000001,  3.984440999999999),
    @(11, 1.387631,           0.5094),
    @(12, 1.179422,           1.035511),
    @(13, 3.411584999999999,  2.3346),
    @(14, 3.319501,           2.484423),
    @(15, 4.815894,           3.980816),
    @(16, 4.317195000000001,  3.482116999999999),
    @(17, 4.835497999999999,  4.00042),
    @(18, 1.394468,           0.379851),
    @(19, 9.651191000000001,  8.816113),
    @(20, 9.657482,           8.822403999999999),
    @(21, 6.884347,           6.049268999999999),
    @(22, 7.215117,           6.380038999999999),
    @(23, 4.060592000000001,  3.225514)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
